$wb = $excel.ActiveWorkbook

# --- Update shared-string-backed cell values on sheets "B1" and "B2" ---
# B1 (formerly showing TxHash index 12/13 placeholders) gets its own unique hashes
$wsB1 = $wb.Worksheets.Item("B1")
$wsB1.Range("A2").Value = "F18639EC56204A7CA80D622B8AC7E97915DD5BE0F4BEAAD7CC77B7F4214CEC34"
$wsB1.Range("A3").Value = "573D1F41A330F78FB529AA265DAE5CBAF4B3A27EC22E78A0A6EDA5972491D38D"

# B2 gets its own unique hashes
$wsB2 = $wb.Worksheets.Item("B2")
$wsB2.Range("A2").Value = "2ABEB91AD183C1FEF536E6C148BA8537DA3C14E59782B18BEAF9E8BBEEC2FC35"
$wsB2.Range("A3").Value = "C47D40FEDD91AB06F879E289CEBBBA5FC26E81A2685562CE7D927F3FE061CCA5"

# --- Update the selection (active cell) on each relevant sheet ---
# A12 previously held the active/selected tab; it keeps its own selection (H28)
# but will no longer be the active sheet once another sheet is activated below.

# B2's selection moves from D7 to M22 (B2 stays a non-active, background tab)
$wsB2.Range("M22").Select()

# B1's selection moves from A2 to K24, and B1 becomes the active sheet/tab
$wsB1.Range("K24").Select()
